# Refresh cryptocurrency price/volume figures (GitHub Actions bot update).
# D = Price column, E = Volume(1h) column. Both are plain text cells in the
# source sheet (thousands separated with "." rather than ",", percentages kept
# as padded strings), so numeric-looking D values must be forced to Text via
# NumberFormat "@" before assignment -- otherwise Excel auto-coerces them into
# real numbers (e.g. "21.00" -> 21, "1.005" -> 1.0049999999999999) and the
# trailing-zero / multi-dot formatting from the source data is lost.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.539.79'
$ws.Range("D3").Value = '1.791.43'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.59'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5347'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3795'
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07503'
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.41'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.114'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.00'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.166'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.370'
$ws.Range("E15").Value = '  +5.44%  '
$ws.Range("D16").Value = '1.794.07'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.04'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06440'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.916'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '28.560.62'
$ws.Range("E23").Value = '  +4.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.20'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.134'
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.41'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.42'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '2.000.45'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.353'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.35'
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.119'
$ws.Range("E31").Value = '  +3.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1017'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.654'
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.660'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2290'
$ws.Range("E35").Value = '  +11.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06555'
$ws.Range("E36").Value = '  +9.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02317'
$ws.Range("E37").Value = '  +2.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.683'
$ws.Range("E38").Value = '  +5.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.070'
$ws.Range("E39").Value = '  +3.00%  '
$ws.Range("E40").Value = '  +1.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6297'
$ws.Range("E41").Value = '  +2.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.207'
$ws.Range("E42").Value = '  +6.60%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.394'
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.61'
$ws.Range("E45").Value = '  +2.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5913'
$ws.Range("E46").Value = '  +2.14%  '
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.61'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.974'
$ws.Range("E49").Value = '  +4.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.156'
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06916'
$ws.Range("E51").Value = '  +2.79%  '
